# Atualização de bases das ligas, do dia: 14-06-2024 às 20:31
# The source feed had four fixture rows associated with the wrong
# match id / odds set (rows swapped pairwise). Fix by swapping the
# row data (columns B:AD -- everything except the running index in
# column A) between each pair of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param(
        [int]$row1,
        [int]$row2
    )
    $rng1 = $ws.Range("B$row1`:AD$row1")
    $rng2 = $ws.Range("B$row2`:AD$row2")

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}

Swap-RowData 14 15
Swap-RowData 19 20
Swap-RowData 25 26
Swap-RowData 105 106
